$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for rows 2-24, per the target diff (row contents were reshuffled).
$rows = @{
  2  = @(44509, 20, 4000, 4000, 4000, 4000)
  3  = @(44280, 55, 4000, 4000, 4000, 4000)
  4  = @(44176, 10, 4000, 4000, 4000, 4000)
  5  = @(44749, 65, 6000, 6000, 6000, 6000)
  6  = @(44504, 55, 4000, 4000, 4000, 4000)
  7  = @(44777, 25, 5000, 5000, 5000, 5000)
  8  = @(44390, 55, 6000, 6000, 6000, 6000)
  9  = @(44498, 40, 4000, 4000, 4000, 4000)
  10 = @(44259, 30, 4000, 4000, 4000, 4000)
  11 = @(44781, 40, 5000, 5000, 5000, 5000)
  12 = @(44497, 20, 4000, 4000, 4000, 4000)
  13 = @(44312, 50, 4000, 4000, 4000, 4000)
  14 = @(44316, 20, 4000, 4000, 4000, 4000)
  15 = @(44680, 20, 5000, 5000, 5000, 5000)
  16 = @(44315, 40, 4000, 4000, 4000, 4000)
  17 = @(44313, 20, 4000, 4000, 4000, 4000)
  18 = @(44508, 30, 4000, 4000, 4000, 4000)
  19 = @(44301, 40, 3000, 3000, 3000, 3000)
  20 = @(44679, 50, 5000, 5000, 5000, 5000)
  21 = @(44649, 20, 5000, 5000, 5000, 5000)
  22 = @(44656, 85, 5000, 5000, 5000, 5000)
  23 = @(44365, 55, 5000, 5000, 5000, 5000)
  24 = @(44291, 35, 4000, 4000, 4000, 4000)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
  $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
  $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
  $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
  $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
  $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio $/Kg
}
